$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C10 (Rule R30, "From" threshold) from 18 to 1
$ws.Range("C10").Value = 1
